# Fixed spelling errors in validation tables.
#
# The workbook used "BioGears" as the name of the physiology engine in a
# series of duplicated label strings (e.g. "BioGears HeartRate" next to
# "HeartRate") and in the scenario narrative text. The engine has since
# been renamed; update all of those labels to "Engine" and update the
# narrative blurb to reference the "Pulse" physiology engine instead of
# "BioGears".

$wb = $excel.ActiveWorkbook

# Rename every "BioGears <Metric>" label (on every worksheet) to
# "Engine <Metric>" - this covers the duplicated column headers on the
# "Heat Stroke Breakdown" sheet (BioGears HeartRate, BioGears
# HeartStrokeVolume, BioGears BloodVolume, BioGears MeanArterialPressure,
# BioGears SystolicArterialPressure, BioGears DiastolicArterialPressure,
# BioGears CardiacOutput, BioGears RespirationRate, BioGears
# OxygenSaturation, BioGears CoreTemperature, BioGears SkinTemperature,
# BioGears SweatRate) as well as the footnote "Note that BioGears tracks
# ionized calcium...".
foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace("BioGears ", "Engine ")
}

# Update the scenario overview narrative on the "Heat Stroke Overview"
# sheet, which referenced "BioGears(R) physiology engine" by name; it now
# calls out the "Pulse" physiology engine instead.
$overview = $wb.Worksheets.Item("Heat Stroke Overview")
$overview.Range("A2").Value = "The heat stroke scenario simulates the body's temperature regulation system. This scenario highlights the ability of the Pulse physiology engine to simulate the energy exchange between the human body and the enviroment."
